$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header rich-text strings ---
$volRng = $ws.Range("A8")
$volRng.Characters(21, 1).Text = "8"
$weekRng = $ws.Range("C9")
$weekRng.Characters(27, 9).Text = "2/19/2024"
$weekRng.Characters(47, 9).Text = "2/25/2024"

# --- Update CompStat table (rows 15-29) ---
# Plain numeric value updates (no type/style change needed)
$ws.Cells.Item(15, 7).Value2 = 2
$ws.Cells.Item(15, 8).Value2 = -50
$ws.Cells.Item(15, 10).Value2 = 2
$ws.Cells.Item(15, 11).Value2 = -50
$ws.Cells.Item(15, 12).Value2 = -75
$ws.Cells.Item(15, 14).Value2 = -92.307692307692
$ws.Cells.Item(16, 3).Value2 = 3
$ws.Cells.Item(16, 4).Value2 = 6
$ws.Cells.Item(16, 5).Value2 = -50
$ws.Cells.Item(16, 6).Value2 = 10
$ws.Cells.Item(16, 7).Value2 = 18
$ws.Cells.Item(16, 8).Value2 = -44.444444444444
$ws.Cells.Item(16, 9).Value2 = 31
$ws.Cells.Item(16, 10).Value2 = 40
$ws.Cells.Item(16, 11).Value2 = -22.5
$ws.Cells.Item(16, 12).Value2 = -27.906976744186
$ws.Cells.Item(16, 13).Value2 = -26.190476190476
$ws.Cells.Item(16, 14).Value2 = -86.695278969957
$ws.Cells.Item(17, 3).Value2 = 10
$ws.Cells.Item(17, 4).Value2 = 4
$ws.Cells.Item(17, 5).Value2 = 150
$ws.Cells.Item(17, 6).Value2 = 23
$ws.Cells.Item(17, 7).Value2 = 16
$ws.Cells.Item(17, 8).Value2 = 43.75
$ws.Cells.Item(17, 9).Value2 = 46
$ws.Cells.Item(17, 10).Value2 = 41
$ws.Cells.Item(17, 11).Value2 = 12.195121951219
$ws.Cells.Item(17, 12).Value2 = 64.285714285714
$ws.Cells.Item(17, 13).Value2 = 64.285714285714
$ws.Cells.Item(17, 14).Value2 = -54.455445544554
$ws.Cells.Item(18, 4).Value2 = 2
$ws.Cells.Item(18, 5).Value2 = -100
$ws.Cells.Item(18, 6).Value2 = 7
$ws.Cells.Item(18, 7).Value2 = 12
$ws.Cells.Item(18, 8).Value2 = -41.666666666666
$ws.Cells.Item(18, 10).Value2 = 21
$ws.Cells.Item(18, 11).Value2 = -33.333333333333
$ws.Cells.Item(18, 12).Value2 = -48.148148148148
$ws.Cells.Item(18, 13).Value2 = -68.181818181818
$ws.Cells.Item(18, 14).Value2 = -95.527156549520
$ws.Cells.Item(19, 3).Value2 = 9
$ws.Cells.Item(19, 4).Value2 = 11
$ws.Cells.Item(19, 5).Value2 = -18.181818181818
$ws.Cells.Item(19, 6).Value2 = 51
$ws.Cells.Item(19, 7).Value2 = 37
$ws.Cells.Item(19, 8).Value2 = 37.837837837837
$ws.Cells.Item(19, 9).Value2 = 94
$ws.Cells.Item(19, 10).Value2 = 70
$ws.Cells.Item(19, 11).Value2 = 34.285714285714
$ws.Cells.Item(19, 12).Value2 = -1.052631578947
$ws.Cells.Item(19, 13).Value2 = 100
$ws.Cells.Item(19, 14).Value2 = -44.378698224852
$ws.Cells.Item(20, 3).Value2 = 6
$ws.Cells.Item(20, 5).Value2 = 50
$ws.Cells.Item(20, 6).Value2 = 25
$ws.Cells.Item(20, 7).Value2 = 16
$ws.Cells.Item(20, 8).Value2 = 56.25
$ws.Cells.Item(20, 9).Value2 = 39
$ws.Cells.Item(20, 10).Value2 = 36
$ws.Cells.Item(20, 11).Value2 = 8.333333333333
$ws.Cells.Item(20, 12).Value2 = -46.575342465753
$ws.Cells.Item(20, 13).Value2 = 254.545454545455
$ws.Cells.Item(20, 14).Value2 = -87.850467289719
$ws.Cells.Item(21, 3).Value2 = 28
$ws.Cells.Item(21, 4).Value2 = 28
$ws.Cells.Item(21, 5).Value2 = 0
$ws.Cells.Item(21, 6).Value2 = 117
$ws.Cells.Item(21, 8).Value2 = 15.841584158415
$ws.Cells.Item(21, 9).Value2 = 225
$ws.Cells.Item(21, 10).Value2 = 210
$ws.Cells.Item(21, 11).Value2 = 7.142857142857
$ws.Cells.Item(21, 12).Value2 = -16.666666666666
$ws.Cells.Item(21, 13).Value2 = 27.840909090909
$ws.Cells.Item(21, 14).Value2 = -80.68669527897
$ws.Cells.Item(22, 6).Value2 = 1
$ws.Cells.Item(22, 7).Value2 = 4
$ws.Cells.Item(22, 8).Value2 = -75
$ws.Cells.Item(22, 10).Value2 = 7
$ws.Cells.Item(22, 11).Value2 = 0
$ws.Cells.Item(22, 12).Value2 = -22.222222222222
$ws.Cells.Item(22, 13).Value2 = 0
$ws.Cells.Item(23, 6).Value2 = 5
$ws.Cells.Item(23, 8).Value2 = 150
$ws.Cells.Item(23, 9).Value2 = 7
$ws.Cells.Item(23, 11).Value2 = 40
$ws.Cells.Item(23, 12).Value2 = 133.333333333333
$ws.Cells.Item(23, 13).Value2 = 16.666666666666
$ws.Cells.Item(24, 3).Value2 = 27
$ws.Cells.Item(24, 4).Value2 = 24
$ws.Cells.Item(24, 5).Value2 = 12.5
$ws.Cells.Item(24, 7).Value2 = 89
$ws.Cells.Item(24, 8).Value2 = 7.865168539325
$ws.Cells.Item(24, 9).Value2 = 190
$ws.Cells.Item(24, 10).Value2 = 196
$ws.Cells.Item(24, 11).Value2 = -3.061224489795
$ws.Cells.Item(24, 12).Value2 = -10.798122065727
$ws.Cells.Item(24, 13).Value2 = 113.483146067416
$ws.Cells.Item(25, 3).Value2 = 11
$ws.Cells.Item(25, 4).Value2 = 11
$ws.Cells.Item(25, 5).Value2 = 0
$ws.Cells.Item(25, 6).Value2 = 49
$ws.Cells.Item(25, 7).Value2 = 40
$ws.Cells.Item(25, 8).Value2 = 22.5
$ws.Cells.Item(25, 9).Value2 = 81
$ws.Cells.Item(25, 10).Value2 = 62
$ws.Cells.Item(25, 11).Value2 = 30.645161290322
$ws.Cells.Item(25, 12).Value2 = 28.571428571428
$ws.Cells.Item(25, 13).Value2 = 5.194805194805
$ws.Cells.Item(26, 7).Value2 = 2
$ws.Cells.Item(26, 8).Value2 = 0
$ws.Cells.Item(26, 10).Value2 = 4
$ws.Cells.Item(26, 11).Value2 = -50
$ws.Cells.Item(26, 12).Value2 = -71.428571428571
$ws.Cells.Item(27, 6).Value2 = 4
$ws.Cells.Item(27, 7).Value2 = 3
$ws.Cells.Item(27, 8).Value2 = 33.333333333333
$ws.Cells.Item(27, 9).Value2 = 7
$ws.Cells.Item(27, 11).Value2 = -12.5
$ws.Cells.Item(27, 12).Value2 = 16.666666666666
$ws.Cells.Item(28, 8).Value2 = -100
$ws.Cells.Item(28, 14).Value2 = -86.842105263157
$ws.Cells.Item(29, 8).Value2 = -100
$ws.Cells.Item(29, 14).Value2 = -86.486486486486

# Type-changing cells: fix style via paste-format from a stable reference cell, then set value
$ws.Cells.Item(15, 4).Value2 = 1
$ws.Cells.Item(16, 3).Copy() | Out-Null
$ws.Cells.Item(15, 4).PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(15, 5).Value2 = -100
$ws.Cells.Item(15, 8).Copy() | Out-Null
$ws.Cells.Item(15, 5).PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(18, 3).NumberFormat = "@"
$ws.Cells.Item(18, 3).Value2 = "0"
$ws.Cells.Item(15, 3).Copy() | Out-Null
$ws.Cells.Item(18, 3).PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(22, 4).Value2 = 2
$ws.Cells.Item(16, 3).Copy() | Out-Null
$ws.Cells.Item(22, 4).PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(22, 5).Value2 = -50
$ws.Cells.Item(15, 8).Copy() | Out-Null
$ws.Cells.Item(22, 5).PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(23, 3).Value2 = 3
$ws.Cells.Item(16, 3).Copy() | Out-Null
$ws.Cells.Item(23, 3).PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value2 = "0"
$ws.Cells.Item(15, 3).Copy() | Out-Null
$ws.Cells.Item(23, 4).PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(23, 5).NumberFormat = "@"
$ws.Cells.Item(23, 5).Value2 = "***.*"
$ws.Cells.Item(15, 3).Copy() | Out-Null
$ws.Cells.Item(23, 5).PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(26, 4).Value2 = 1
$ws.Cells.Item(16, 3).Copy() | Out-Null
$ws.Cells.Item(26, 4).PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(26, 5).Value2 = -100
$ws.Cells.Item(15, 8).Copy() | Out-Null
$ws.Cells.Item(26, 5).PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(27, 3).Value2 = 1
$ws.Cells.Item(16, 3).Copy() | Out-Null
$ws.Cells.Item(27, 3).PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value2 = "0"
$ws.Cells.Item(15, 3).Copy() | Out-Null
$ws.Cells.Item(27, 4).PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(27, 5).NumberFormat = "@"
$ws.Cells.Item(27, 5).Value2 = "***.*"
$ws.Cells.Item(15, 3).Copy() | Out-Null
$ws.Cells.Item(27, 5).PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(28, 6).NumberFormat = "@"
$ws.Cells.Item(28, 6).Value2 = "0"
$ws.Cells.Item(15, 3).Copy() | Out-Null
$ws.Cells.Item(28, 6).PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(29, 6).NumberFormat = "@"
$ws.Cells.Item(29, 6).Value2 = "0"
$ws.Cells.Item(15, 3).Copy() | Out-Null
$ws.Cells.Item(29, 6).PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0
